$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.659.02'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.601.72'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('E4').Value = '  -0.11%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '212.14'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.45%  '
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('E7').Value = '  -0.09%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '29.05'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +8.49%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.256'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +2.78%  '
$ws.Range('E10').Value = '  +1.36%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0906'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('D12').Value = '1.831.62'
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('D13').Value = '1.598.55'
$ws.Range('E13').Value = '  -0.74%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.554'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +3.23%  '
$ws.Range('D15').Value = '29.690.69'
$ws.Range('E15').Value = '  +0.57%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '3.78'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.81%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '64.15'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +1.14%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '241.45'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +0.54%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '8.10'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +6.75%  '
$ws.Range('D20').Value = '0.0₃0701'
$ws.Range('E20').Value = '  +1.25%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.10%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.03'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.13%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '9.49'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +3.33%  '
$ws.Range('E24').Value = '  +2.12%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '156.47'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.95%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '15.49'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +1.50%  '
$ws.Range('E27').Value = '  +1.03%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '6.47'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +1.80%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  +2.10%  '
$ws.Range('E31').Value = '  -0.84%  '
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('E33').Value = '  +2.15%  '
$ws.Range('D34').Value = '1.422.03'
$ws.Range('E34').Value = '  -0.71%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.57'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +4.17%  '
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('E37').Value = '  +0.63%  '
$ws.Range('E38').Value = '  +0.63%  '
$ws.Range('E39').Value = '  +2.34%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.546'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +2.76%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '55.73'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +5.57%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.0496'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +6.67%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.818'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +3.11%  '
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('E45').Value = '  -0.17%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '67.36'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +2.33%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.994'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +18.87%  '
$ws.Range('E48').Value = '  +2.26%  '
$ws.Range('D49').Value = '1.741.13'
$ws.Range('E49').Value = '  -0.26%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '86.66'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('E51').Value = '  -1.48%  '
